# The workbook uses "##" as an internal separator between a comment/value
# and a trailing date in many cells (e.g. "0551-01993-25##01.07.2025").
# This edit replaces every occurrence of "##" with a single space " "
# across all used cells of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange
[void]$ur.Replace("##", " ")
